# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-19 09:29:12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Missing Sessions metric values (K7/L7 and K8/L8 area) ---
$ws.Range("L7").Value = 51
$ws.Range("L8").Value = 18

# --- 2. "Recorded By" column: swap "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$gRows = @(8,9,10,12,14,15,17,18,23,34,35,36,38,40,41,43,44,49,60,61,62,64,66,67,69,70,75,86,87,88,90,92,93,95,96,101,112,113,114,116,118,119,121,122,127,138,139,140,142,144,145,147,148,153,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($r in $gRows) {
    $ws.Range("G" + $r).Value = "System, dnasr281@gmail.com"
}

# --- 3. Per-row Missing/Excused counters (rows 15-20, columns P & Q) ---
$ws.Range("P15").Value = 5
$ws.Range("Q15").Value = 0
$ws.Range("P16").Value = 4
$ws.Range("Q16").Value = 0
$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 0
$ws.Range("P18").Value = 4
$ws.Range("Q18").Value = 0
$ws.Range("P19").Value = 4
$ws.Range("Q19").Value = 0
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = 0

# --- 4. Upcoming-session summary rows: re-style from "Pending" (yellow) to "Not Recorded" (red)
#        and update the status text, using row 3 (an existing "Not Recorded" row) as the format source ---
$fmtSource = $ws.Range("A3:I3")
$fmtSource.Copy()
$notRecordedRows = @(27,53,79,105,131,157)
foreach ($r in $notRecordedRows) {
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $dst.PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = "Not Recorded"
}

Write-Output "done"
